# Update cryptos list (cell values) per commit "Updated cryptos list on Thu Nov  9 10:49:16 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.822.60"
$ws.Range("E2").Value = "  +4.08%  "
$ws.Range("D3").Value = "1.914.74"
$ws.Range("E3").Value = "  +1.55%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'250.56"
$ws.Range("E5").Value = "  +1.55%  "
$ws.Range("D6").Value = "'0.702"
$ws.Range("E6").Value = "  +0.50%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'46.85"
$ws.Range("E8").Value = "  +8.23%  "
$ws.Range("E9").Value = "  +4.42%  "
$ws.Range("D10").Value = "'58.16"
$ws.Range("E10").Value = "  +8.79%  "
$ws.Range("D11").Value = "'0.0762"
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("E12").Value = "  +2.14%  "
$ws.Range("D13").Value = "'14.71"
$ws.Range("E13").Value = "  +8.44%  "
$ws.Range("D14").Value = "'0.819"
$ws.Range("E14").Value = "  +5.53%  "
$ws.Range("D15").Value = "2.191.59"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("E16").Value = "  +3.44%  "
$ws.Range("D17").Value = "1.917.91"
$ws.Range("E17").Value = "  +1.83%  "
$ws.Range("D18").Value = "37.260.80"
$ws.Range("E18").Value = "  +5.40%  "
$ws.Range("D19").Value = "'74.78"
$ws.Range("E19").Value = "  +1.33%  "
$ws.Range("D20").Value = "0.0₃0859"
$ws.Range("E20").Value = "  +3.54%  "
$ws.Range("D21").Value = "'13.64"
$ws.Range("E21").Value = "  +5.96%  "
$ws.Range("D22").Value = "'251.55"
$ws.Range("E22").Value = "  +2.48%  "
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("D24").Value = "'2.61"
$ws.Range("E24").Value = "  -0.48%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("E26").Value = "  +3.56%  "
$ws.Range("D27").Value = "'167.98"
$ws.Range("E27").Value = "  +1.80%  "
$ws.Range("E28").Value = "  +1.53%  "
$ws.Range("D29").Value = "'18.73"
$ws.Range("E29").Value = "  +2.04%  "
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("D31").Value = "'4.62"
$ws.Range("E31").Value = "  +7.20%  "
$ws.Range("D32").Value = "'0.0621"
$ws.Range("E32").Value = "  +4.04%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'4.35"
$ws.Range("E33").Value = "  +3.54%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").Value = "'0.0896"
$ws.Range("E34").Value = "  +20.92%  "
$ws.Range("D35").Value = "'1.91"
$ws.Range("E35").Value = "  +1.81%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("B37").Value = "Gas"
$ws.Range("C37").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D37").Value = "'19.21"
$ws.Range("E37").Value = "  +58.33%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.51"
$ws.Range("E38").Value = "  +2.72%  "
$ws.Range("D39").Value = "'0.876"
$ws.Range("E39").Value = "  +1.74%  "
$ws.Range("E40").Value = "  +2.63%  "
$ws.Range("D41").Value = "'105.69"
$ws.Range("E41").Value = "  +8.17%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "'18.02"
$ws.Range("E42").Value = "  +4.07%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0228"
$ws.Range("E43").Value = "  +4.41%  "
$ws.Range("D44").Value = "'2.88"
$ws.Range("E44").Value = "  +20.17%  "
$ws.Range("E45").Value = "  +2.13%  "
$ws.Range("D46").Value = "1.348.94"
$ws.Range("E46").Value = "  +2.67%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "'0.0818"
$ws.Range("E48").Value = "  +1.43%  "
$ws.Range("D49").Value = "'2.83"
$ws.Range("E49").Value = "  +3.50%  "
$ws.Range("E50").Value = "  +2.12%  "
$ws.Range("D51").Value = "'43.13"
$ws.Range("E51").Value = "  +1.14%  "
